# 8 August 2017 barkha
# Add the "successfully saved" message strings (Message_Text1/Message_Text2 +
# their localized values) into row 3 of the Portuguese sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portuguese")

$ws.Range("E3").Value = "Message_Text1"
# Leading apostrophe forces this text (which starts with a hyphen) to be
# stored as literal text with the quote-prefix cell style - matches the
# workbook's existing convention used on H2.
$ws.Range("F3").Value = "'-  Amenidade/serviço"
$ws.Range("G3").Value = "Message_Text2"
$ws.Range("H3").Value = " salvo com sucesso."

# Widen columns E, G and H to fit the new content.
$ws.Columns.Item(5).ColumnWidth = 13.833333333333334
$ws.Columns.Item(7).ColumnWidth = 13.833333333333334
$ws.Columns.Item(8).ColumnWidth = 17.5

# Move the active selection to H6, as left by the author.
$ws.Range("H6").Select()
